$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.646.92'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '3.329.64'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.78%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -1.83%  '
$ws.Range('D9').Value = '3.324.87'
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.178'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.576'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('E12').Value = '  -3.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000270'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '663.98'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.28%  '
$ws.Range('D15').Value = '3.871.78'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.41'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('D17').Value = '67.820.52'
$ws.Range('E17').Value = '  -1.09%  '
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').Value = '3.330.65'
$ws.Range('E19').Value = '  -1.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.01%  '
$ws.Range('E21').Value = '  -1.51%  '
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.42'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.09%  '
$ws.Range('E27').Value = '  -6.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.45'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.71%  '
$ws.Range('E31').Value = '  -2.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '591.79'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.95'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.48%  '
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('D36').Value = '3.715.13'
$ws.Range('E36').Value = '  -6.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.94%  '
$ws.Range('E38').Value = '  -12.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.131'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('E41').Value = '  -6.06%  '
$ws.Range('E42').Value = '  -6.68%  '
$ws.Range('E43').Value = '  -3.06%  '
$ws.Range('D44').Value = '0.0₃0664'
$ws.Range('E44').Value = '  -5.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.28%  '
$ws.Range('E46').Value = '  -3.96%  '
$ws.Range('E47').Value = '  -0.60%  '
$ws.Range('E48').Value = '  -1.74%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '127.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.26%  '
